$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.052.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.32%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.834.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.66%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6258'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.69%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07569'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.75%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2924'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.41%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.86%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07727'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.830.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.956'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6636'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001011'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +18.01%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.70'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.041'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.13%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '28.998.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.09%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '226.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.55%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.02%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.183'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.92%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.56%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.497'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.04%  '

$ws.Range('E26').Value = '  +0.74%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.44%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.491'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.68%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.098'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.79%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.017'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.62%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.194'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.50%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05241'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.51%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.841'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.08%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7336'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.30%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.138'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.95%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.699'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.99%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.236.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.93%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.763'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.89%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01781'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.46%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.321'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.11%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8974'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.52%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.29%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.972.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.37%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000124'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.46%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5114'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.44%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4035'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.99%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.852'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.87%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05749'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.668'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.48%  '
